# Bugfix the naive forecaster component: the "date" column (column A) was
# storing raw Excel date serials (formatted with a custom date/time number
# format). It should instead hold plain "YYYYQn" quarter-label text, matching
# the header style used by the rest of the table (centered, bold, bordered).
#
# For every data row we:
#   1. Read the existing date serial value out of column A.
#   2. Turn it into the corresponding "YYYYQn" label.
#   3. Write that label back as text.
#   4. Copy the header cell's formatting (font/border/alignment) onto the
#      cell so it no longer uses the old date-specific number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Row 1 is the header ("value" / "date"); data starts on row 2.
$headerRow = $firstRow

# Copy the header's cell format once; we'll paste it onto every data cell in
# column A so they pick up the same style the header already uses (centered,
# bold font, thin border) instead of the old date number format.
$ws.Cells.Item($headerRow, 1).Copy()

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)

    $serial = $cell.Value2
    $d = [DateTime]::FromOADate($serial)
    $quarter = [Math]::Floor(($d.Month - 1) / 3) + 1
    $label = "" + $d.Year + "Q" + $quarter

    $cell.Value = $label
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}
